$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 211
$ws.Range("I2").Value = 438
$ws.Range("J2").Value = 2029
$ws.Range("K2").Value = 11
$ws.Range("L2").Value = 606
$ws.Range("M2").Value = 25
$ws.Range("N2").Value = 373
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 34
$ws.Range("S2").Value = 207
$ws.Range("T2").Value = 357
$ws.Range("U2").Value = 22
$ws.Range("V2").Value = 3251
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 3284
$ws.Range("Z2").Value = 52
$ws.Range("AA2").Value = 17
